# mbr_template.pptx: swap the generic placeholder title slide text for the
# {{MBR_MONTH}} token, and replace the "KPIs (Placeholder Tokens)" title +
# single content placeholder on slide 2 with four free-standing text boxes
# carrying the structured KPI / table / chart tokens.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1: subtitle placeholder text -> {{MBR_MONTH}}
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.Name -eq "Subtitle 2") {
        $shp.TextFrame.TextRange.Text = "{{MBR_MONTH}}"
    }
}

# ---------------------------------------------------------------------
# Slide 2: drop the inherited Title/Content placeholders (each requires
# two Delete() calls since the layout immediately re-seeds a fresh,
# empty instance of a "last of its kind" placeholder) and rebuild the
# slide from plain text boxes.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

while ($s2.Shapes.Count -gt 0) {
    $s2.Shapes.Item(1).Delete()
}
while ($s2.Shapes.Count -gt 0) {
    $s2.Shapes.Item(1).Delete()
}

# -- TextBox 1: "KPIs & Insights" heading -------------------------------
$tb1 = $s2.Shapes.AddTextbox(1, 43.2, 21.6, 864.0, 43.2)
$tb1.Name = "TextBox 1"
$tb1.TextFrame.TextRange.Text = "KPIs & Insights"
$tb1.TextFrame.WordWrap = 0
$tb1.TextFrame.AutoSize = 1
$tb1.TextFrame.TextRange.Font.Size = 28
$tb1.Left = 43.2
$tb1.Top = 21.6
$tb1.Width = 864.0
$tb1.Height = 43.2
$tb1.Fill.Visible = 0

# -- TextBox 2: coverage note / invoice count / net / gross -------------
$tb2 = $s2.Shapes.AddTextbox(1, 43.2, 86.4, 432.0, 100.8)
$tb2.Name = "TextBox 2"
$tb2.TextFrame.TextRange.Text = "{{COVERAGE_NOTE}}`rRechnungen: {{INVOICE_COUNT}}`rNetto: {{TOTAL_NET}}`rBrutto: {{TOTAL_GROSS}}"
$tb2.TextFrame.WordWrap = 0
$tb2.TextFrame.AutoSize = 1
$tb2.Left = 43.2
$tb2.Top = 86.4
$tb2.Width = 432.0
$tb2.Height = 100.8
$tb2.Fill.Visible = 0

# -- TextBox 3: top suppliers table placeholder --------------------------
$tb3 = $s2.Shapes.AddTextbox(1, 43.2, 208.8, 432.0, 273.6)
$tb3.Name = "TextBox 3"
$tb3.TextFrame.TextRange.Text = "{{TOP_SUPPLIERS_TABLE}}"
$tb3.TextFrame.WordWrap = 0
$tb3.TextFrame.AutoSize = 1
$tb3.Left = 43.2
$tb3.Top = 208.8
$tb3.Width = 432.0
$tb3.Height = 273.6
$tb3.Fill.Visible = 0

# -- TextBox 4: budget chart placeholder ---------------------------------
$tb4 = $s2.Shapes.AddTextbox(1, 504.0, 208.8, 403.2, 273.6)
$tb4.Name = "TextBox 4"
$tb4.TextFrame.TextRange.Text = "{{BUDGET_CHART}}"
$tb4.TextFrame.WordWrap = 0
$tb4.TextFrame.AutoSize = 1
$tb4.Left = 504.0
$tb4.Top = 208.8
$tb4.Width = 403.2
$tb4.Height = 273.6
$tb4.Fill.Visible = 0
